$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H59").Value = 4000
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").Value = $null
$ws.Range("H70").Value = 7500
$ws.Range("J70").Value = 7500
$ws.Range("L70").Value = 22500
$ws.Range("N70").Value = -23040
$ws.Range("H73").Value = 7500
$ws.Range("J73").Value = 7500
$ws.Range("L73").Value = 22500
$ws.Range("N73").Value = -24372
$ws.Range("H137").Value = 2001
$ws.Range("I137").Value = 2001
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6003
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3453
$ws.Range("N137").Value = $null

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2000
$ws.Range("I63").Value = 2000
$ws.Range("K63").Value = 2000
$ws.Range("M63").Value = -1314
$ws.Range("H66").Value = 2000
$ws.Range("I66").Value = 2000
$ws.Range("K66").Value = 10000
$ws.Range("M66").Value = -6568
$ws.Range("H110").Value = 2949.2
$ws.Range("I110").Value = 3186.625
$ws.Range("J110").Value = 1999.5
$ws.Range("K110").Value = 3186.625
$ws.Range("L110").Value = 1999.5
$ws.Range("M110").Value = -1141.625
$ws.Range("N110").Value = -6089.5

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 947.5
$ws.Range("I20").Value = 884.875
$ws.Range("K20").Value = 884.875
$ws.Range("M20").Value = -637.875

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 28329.666
$ws.Range("J50").Value = 29994.5
$ws.Range("L50").Value = 29994.5
$ws.Range("N50").Value = -31244.5
$ws.Range("H51").Value = 14549.5
$ws.Range("J51").Value = 20099
$ws.Range("L51").Value = 20099
$ws.Range("N51").Value = -21571
$ws.Range("H58").Value = 4023.6
$ws.Range("I58").Value = 4529.5
$ws.Range("J58").Value = 2000
$ws.Range("K58").Value = 4529.5
$ws.Range("L58").Value = 2000
$ws.Range("M58").Value = -4326.5
$ws.Range("N58").Value = -2406
$ws.Range("H60").Value = 8250
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").Value = $null
$ws.Range("H61").Value = 14549.5
$ws.Range("J61").Value = 20099
$ws.Range("L61").Value = 20099
$ws.Range("N61").Value = -20795
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").Value = $null
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").Value = $null
$ws.Range("H99").Value = 3119.5
$ws.Range("J99").Value = 3082.5
$ws.Range("L99").Value = 3082.5
$ws.Range("N99").Value = -6078.5
$ws.Range("H105").Value = 5894.9
$ws.Range("I105").Value = 2199.6
$ws.Range("J105").Value = 9590.200000000001
$ws.Range("K105").Value = 2199.6
$ws.Range("L105").Value = 9590.200000000001
$ws.Range("M105").Value = -452.5999999999999
$ws.Range("N105").Value = -13084.2
$ws.Range("H122").Value = 963.6667
$ws.Range("I122").Value = 773.6
$ws.Range("K122").Value = 2320.8
$ws.Range("M122").Value = 129.1999999999998
$ws.Range("H126").Value = 3119.5
$ws.Range("J126").Value = 3082.5
$ws.Range("L126").Value = 9247.5
$ws.Range("N126").Value = -14187.5
$ws.Range("H136").Value = 4023.6
$ws.Range("I136").Value = 4529.5
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 13588.5
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -11038.5
$ws.Range("N136").Value = -11100

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 2706.5
$ws.Range("J10").Value = 500
$ws.Range("L10").Value = 1500
$ws.Range("N10").Value = -1778
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = $null
$ws.Range("H59").Value = 815.5
$ws.Range("I59").Value = 131
$ws.Range("K59").Value = 393
$ws.Range("M59").Value = 147
$ws.Range("H129").Value = 1548.1111
$ws.Range("J129").Value = 2437.25
$ws.Range("L129").Value = 7311.75
$ws.Range("N129").Value = -17311.75
$ws.Range("H139").Value = 4219.3
$ws.Range("J139").Value = 4000
$ws.Range("L139").Value = 12000
$ws.Range("N139").Value = -22280
$ws.Range("H140").Value = 3606.0557
$ws.Range("I140").Value = 681.8
$ws.Range("K140").Value = 2045.4
$ws.Range("M140").Value = 3134.6

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2328
$ws.Range("I113").Value = 2104
$ws.Range("K113").Value = 2104
$ws.Range("M113").Value = 66

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 1002
$ws.Range("I62").Value = 1002
$ws.Range("K62").Value = 1002
$ws.Range("M62").Value = -378
$ws.Range("H65").Value = 1002
$ws.Range("I65").Value = 1002
$ws.Range("K65").Value = 5010
$ws.Range("M65").Value = -1890
$ws.Range("H81").Value = 3498
$ws.Range("I81").Value = 3747.5
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 7495
$ws.Range("L81").Value = 5000
$ws.Range("M81").Value = -6434
$ws.Range("N81").Value = -7122
$ws.Range("H84").Value = 3498
$ws.Range("I84").Value = 3747.5
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 37475
$ws.Range("L84").Value = 25000
$ws.Range("M84").Value = -32171
$ws.Range("N84").Value = -35608
$ws.Range("H122").Value = 4010.6667
$ws.Range("J122").Value = 3333
$ws.Range("L122").Value = 9999
$ws.Range("N122").Value = -14899
$ws.Range("H126").Value = 1368.6666
$ws.Range("I126").Value = 1497.5
$ws.Range("J126").Value = 1111
$ws.Range("K126").Value = 4492.5
$ws.Range("L126").Value = 3333
$ws.Range("M126").Value = -2022.5
$ws.Range("N126").Value = -8273
$ws.Range("H136").Value = 884.375
$ws.Range("I136").Value = 867.8570999999999
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 2603.5713
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -53.57129999999961
$ws.Range("N136").Value = -8100
